{"js": "// Update the East Asian / complex-script fonts recorded on the document's\n// paragraph styles (vignettes/docx/ex-tbl.docx style sheet update).\n//\n//   - \"Normal\" and \"Heading\" styles: eastAsia font DejaVu Sans -> Tahoma\n//   - \"List\", \"Caption\", \"Index\" styles: add an explicit complex-script\n//     (cs) font of DejaVu Sans (previously inherited / unset)\n\nconst styles = context.document.getStyles();\n\nconst normal = styles.getByNameOrNullObject(\"Normal\");\nconst heading = styles.getByNameOrNullObject(\"Heading\");\nconst list = styles.getByNameOrNullObject(\"List\");\nconst caption = styles.getByNameOrNullObject(\"Caption\");\nconst index = styles.getByNameOrNullObject(\"Index\");\nawait context.sync();\n\nnormal.font.nameFarEast = \"Tahoma\";\nheading.font.nameFarEast = \"Tahoma\";\nlist.font.nameBidirectional = \"DejaVu Sans\";\ncaption.font.nameBidirectional = \"DejaVu Sans\";\nindex.font.nameBidirectional = \"DejaVu Sans\";\n\nawait context.sync();\n", "ps1": "# Update the East Asian / complex-script fonts recorded on the document's\n# paragraph styles (vignettes/docx/ex-tbl.docx style sheet update).\n#\n#   - \"Normal\" and \"Heading\" styles: eastAsia font DejaVu Sans -> Tahoma\n#   - \"List\", \"Caption\", \"Index\" styles: add an explicit complex-script\n#     (cs) font of DejaVu Sans (previously inherited / unset)\n\n$d = $word.ActiveDocument\n\n$normal = $d.Styles.Item(\"Normal\")\n$normal.Font.NameFarEast = \"Tahoma\"\n\n$heading = $d.Styles.Item(\"Heading\")\n$heading.Font.NameFarEast = \"Tahoma\"\n\n$list = $d.Styles.Item(\"List\")\n$list.Font.NameBi = \"DejaVu Sans\"\n\n$caption = $d.Styles.Item(\"Caption\")\n$caption.Font.NameBi = \"DejaVu Sans\"\n\n$index = $d.Styles.Item(\"Index\")\n$index.Font.NameBi = \"DejaVu Sans\"\n"}
